# Apply the "Update countries & provincias Spain" edit:
#  - refresh the "Datos actualizados" timestamp
#  - re-rank India/Peru, Egipto/Sudafrica/Chequia/Noruega and Jordania/Gabon
#    (each now sits in the row matching its updated case count)
#  - update the daily case-count figures for the affected rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 18:34"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1299770
$ws.Cells.Item(4, 3).Value = 7147
$ws.Cells.Item(4, 4).Value = 219482
$ws.Cells.Item(4, 5).Value = 1002729
$ws.Cells.Item(4, 6).Value = 17020
$ws.Cells.Item(4, 7).Value = 631
$ws.Cells.Item(4, 8).Value = 77559

# Row 6: Italia
$ws.Cells.Item(6, 2).Value = 217185
$ws.Cells.Item(6, 3).Value = 1327
$ws.Cells.Item(6, 4).Value = 99023
$ws.Cells.Item(6, 5).Value = 87961
$ws.Cells.Item(6, 6).Value = 1168
$ws.Cells.Item(6, 7).Value = 243
$ws.Cells.Item(6, 8).Value = 30201

# Row 7: Reino Unido
$ws.Cells.Item(7, 2).Value = 211364
$ws.Cells.Item(7, 3).Value = 4649
$ws.Cells.Item(7, 5).Value = 179779
$ws.Cells.Item(7, 7).Value = 626
$ws.Cells.Item(7, 8).Value = 31241

# Row 11: Brasil
$ws.Cells.Item(11, 2).Value = 138121
$ws.Cells.Item(11, 3).Value = 2428
$ws.Cells.Item(11, 5).Value = 73381
$ws.Cells.Item(11, 7).Value = 202
$ws.Cells.Item(11, 8).Value = 9390

# Row 15: Canada
$ws.Cells.Item(15, 2).Value = 65400
$ws.Cells.Item(15, 3).Value = 478
$ws.Cells.Item(15, 4).Value = 29682
$ws.Cells.Item(15, 5).Value = 31245
$ws.Cells.Item(15, 7).Value = 65
$ws.Cells.Item(15, 8).Value = 4473

# Row 16: India
$ws.Cells.Item(16, 1).Value = "India"
$ws.Cells.Item(16, 2).Value = 59205
$ws.Cells.Item(16, 3).Value = 2854
$ws.Cells.Item(16, 4).Value = 17698
$ws.Cells.Item(16, 5).Value = 39527
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 91
$ws.Cells.Item(16, 8).Value = 1980

# Row 17: Peru
$ws.Cells.Item(17, 1).Value = "Peru"
$ws.Cells.Item(17, 2).Value = 58526
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 18388
$ws.Cells.Item(17, 5).Value = 38511
$ws.Cells.Item(17, 6).Value = 722
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 1627

# Row 29: Singapur
$ws.Cells.Item(29, 4).Value = 2040
$ws.Cells.Item(29, 5).Value = 19647
$ws.Cells.Item(29, 6).Value = 22

# Row 36: Polonia
$ws.Cells.Item(36, 2).Value = 15366
$ws.Cells.Item(36, 3).Value = 319
$ws.Cells.Item(36, 5).Value = 9406
$ws.Cells.Item(36, 7).Value = 21
$ws.Cells.Item(36, 8).Value = 776

# Row 37: Rumania
$ws.Cells.Item(37, 5).Value = 7465
$ws.Cells.Item(37, 7).Value = 35
$ws.Cells.Item(37, 8).Value = 923

# Row 47: Egipto
$ws.Cells.Item(47, 1).Value = "Egipto"
$ws.Cells.Item(47, 2).Value = 8476
$ws.Cells.Item(47, 3).Value = 495
$ws.Cells.Item(47, 4).Value = 1945
$ws.Cells.Item(47, 5).Value = 6028
$ws.Cells.Item(47, 6).Value = 41
$ws.Cells.Item(47, 7).Value = 21
$ws.Cells.Item(47, 8).Value = 503

# Row 48: Sudafrica
$ws.Cells.Item(48, 1).Value = "Sudafrica"
$ws.Cells.Item(48, 2).Value = 8232
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 3153
$ws.Cells.Item(48, 5).Value = 4918
$ws.Cells.Item(48, 6).Value = 36
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 161

# Row 49: Chequia
$ws.Cells.Item(49, 1).Value = "Chequia"
$ws.Cells.Item(49, 2).Value = 8065
$ws.Cells.Item(49, 3).Value = 34
$ws.Cells.Item(49, 4).Value = 4408
$ws.Cells.Item(49, 5).Value = 3386
$ws.Cells.Item(49, 6).Value = 258
$ws.Cells.Item(49, 8).Value = 271

# Row 50: Noruega
$ws.Cells.Item(50, 1).Value = "Noruega"
$ws.Cells.Item(50, 2).Value = 8055
$ws.Cells.Item(50, 3).Value = 21
$ws.Cells.Item(50, 4).Value = 32
$ws.Cells.Item(50, 5).Value = 7805
$ws.Cells.Item(50, 6).Value = 27
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 218

# Row 56: Marruecos
$ws.Cells.Item(56, 2).Value = 5711
$ws.Cells.Item(56, 3).Value = 163
$ws.Cells.Item(56, 4).Value = 2324
$ws.Cells.Item(56, 5).Value = 3201
$ws.Cells.Item(56, 7).Value = 3
$ws.Cells.Item(56, 8).Value = 186

# Row 62: Luxemburgo
$ws.Cells.Item(62, 2).Value = 3871
$ws.Cells.Item(62, 3).Value = 12
$ws.Cells.Item(62, 4).Value = 3526
$ws.Cells.Item(62, 5).Value = 245
$ws.Cells.Item(62, 6).Value = 16

# Row 119: Jordania
$ws.Cells.Item(119, 1).Value = "Jordania"
$ws.Cells.Item(119, 2).Value = 508
$ws.Cells.Item(119, 3).Value = 14
$ws.Cells.Item(119, 4).Value = 385
$ws.Cells.Item(119, 5).Value = 114
$ws.Cells.Item(119, 6).Value = 5
$ws.Cells.Item(119, 8).Value = 9

# Row 120: Gabon
$ws.Cells.Item(120, 1).Value = "Gabon"
$ws.Cells.Item(120, 2).Value = 504
$ws.Cells.Item(120, 4).Value = 110
$ws.Cells.Item(120, 5).Value = 386
$ws.Cells.Item(120, 6).Value = 1
$ws.Cells.Item(120, 8).Value = 8
